$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8623912930488586
$ws.Range("B1").Value = 2.860342025756836
$ws.Range("C1").Value = 8.763140678405762
$ws.Range("D1").Value = 2.025200843811035
$ws.Range("E1").Value = 1.145461320877075
